$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Explanation of Back Propagation"

$ws.Range("A4").Select()
